$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Column D holds numeric-looking text (e.g. "60.694.85", "5.60") that must
# stay literal text rather than being auto-coerced to a number by Excel.
# Format the whole D:E data range as Text first, write the values, then
# restore the original (default) style so cell formatting is unaffected.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "60.694.85"
$ws.Range("E2").Value = "  -5.34%  "
$ws.Range("D3").Value = "3.271.69"
$ws.Range("E3").Value = "  -5.97%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "562.02"
$ws.Range("E5").Value = "  -3.87%  "
$ws.Range("D6").Value = "126.47"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.272.89"
$ws.Range("E8").Value = "  -5.91%  "
$ws.Range("D10").Value = "7.28"
$ws.Range("E10").Value = "  -4.59%  "
$ws.Range("E11").Value = "  -4.63%  "
$ws.Range("E12").Value = "  -4.16%  "
$ws.Range("D13").Value = "3.842.76"
$ws.Range("E13").Value = "  -5.55%  "
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "3.278.52"
$ws.Range("E15").Value = "  -5.80%  "
$ws.Range("E16").Value = "  -6.09%  "
$ws.Range("D17").Value = "60.886.07"
$ws.Range("E17").Value = "  -5.04%  "
$ws.Range("D18").Value = "24.16"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "5.60"
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "13.19"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").Value = "8.93"
$ws.Range("E21").Value = "  -10.36%  "
$ws.Range("D22").Value = "348.61"
$ws.Range("E22").Value = "  -9.35%  "
$ws.Range("D23").Value = "0.551"
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "3.411.96"
$ws.Range("E25").Value = "  -5.68%  "
$ws.Range("D26").Value = "69.19"
$ws.Range("E26").Value = "  -7.45%  "
$ws.Range("E27").Value = "  -4.60%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "7.15"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "1.42"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("E32").Value = "  -6.12%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "0.149"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").Value = "3.308.20"
$ws.Range("E35").Value = "  -5.68%  "
$ws.Range("D36").Value = "22.51"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("D37").Value = "5.21"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "6.76"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "159.70"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D41").Value = "0.0747"
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "40.96"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "0.736"
$ws.Range("E45").Value = "  -7.81%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("E47").Value = "  -4.42%  "
$ws.Range("D48").Value = "22.40"
$ws.Range("E48").Value = "  -5.08%  "
$ws.Range("D49").Value = "6.64"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").Value = "0.854"
$ws.Range("E50").Value = "  -5.51%  "
$ws.Range("D51").Value = "21.20"
$ws.Range("E51").Value = "  +3.76%  "

$dataRange.Style = "Normal"
